# Update "想去人数" (F column) values on the "展览" sheet
$wb = $excel.ActiveWorkbook

$wsZhanLan = $wb.Worksheets.Item("展览")
$wsZhanLan.Range("F2").Value = 78
$wsZhanLan.Range("F3").Value = 405
$wsZhanLan.Range("F4").Value = 2626
$wsZhanLan.Range("F5").Value = 455
$wsZhanLan.Range("F7").Value = 24
$wsZhanLan.Range("F8").Value = 258
$wsZhanLan.Range("F9").Value = 14257
$wsZhanLan.Range("F10").Value = 142
$wsZhanLan.Range("F11").Value = 111
$wsZhanLan.Range("F12").Value = 5717
$wsZhanLan.Range("F18").Value = 7
$wsZhanLan.Range("F21").Value = 779
$wsZhanLan.Range("F22").Value = 2922
$wsZhanLan.Range("F23").Value = 51
$wsZhanLan.Range("F24").Value = 10504
$wsZhanLan.Range("F26").Value = 51
$wsZhanLan.Range("F27").Value = 74
$wsZhanLan.Range("F28").Value = 3727

# Update "想去人数" (F column) values on the "全部类型" sheet
$wsQuanBu = $wb.Worksheets.Item("全部类型")
$wsQuanBu.Range("F2").Value = 78
$wsQuanBu.Range("F3").Value = 405
$wsQuanBu.Range("F5").Value = 2649
$wsQuanBu.Range("F6").Value = 455
$wsQuanBu.Range("F8").Value = 24
$wsQuanBu.Range("F9").Value = 258
$wsQuanBu.Range("F10").Value = 14257
$wsQuanBu.Range("F11").Value = 142
$wsQuanBu.Range("F12").Value = 111
$wsQuanBu.Range("F13").Value = 5717
$wsQuanBu.Range("F19").Value = 7
$wsQuanBu.Range("F22").Value = 779
$wsQuanBu.Range("F23").Value = 2922
$wsQuanBu.Range("F24").Value = 51
$wsQuanBu.Range("F26").Value = 10504
$wsQuanBu.Range("F28").Value = 51
$wsQuanBu.Range("F29").Value = 74
$wsQuanBu.Range("F30").Value = 3727
